# Append one new data row (row 98) to the bottom of the daily log table on
# Sheet1, mirroring the existing rows (A: date text, B: weekday text,
# C: hour number, D: ranking number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98

# Column A holds date-looking text (e.g. "2025/10/13") stored as literal text
# in every existing row, not as a real Excel date serial. Assigning the
# string directly would get auto-parsed into a date value, so force it to be
# treated as text first, then drop the resulting "Text" number-format style
# back to Normal so the new cell matches the unstyled data cells above it.
$ws.Range("A" + $row).Value = "'2025/10/13"
$ws.Range("A" + $row).Style = "Normal"

$ws.Range("B" + $row).Value = "月"
$ws.Range("C" + $row).Value = 5
$ws.Range("D" + $row).Value = 201
